$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new column E: per-row category image asset path (dynamic image
# category loading). Rows 1-5 = Airport, 6-10 = Getting Around, 11-15 = Restaurant.
$airport = "file:///android_asset/airport.png"
$gettingAround = "file:///android_asset/getting_around.png"
$restaurant = "file:///android_asset/restaurant.png"

$ws.Range("E1").Value = $airport
$ws.Range("E2").Value = $airport
$ws.Range("E3").Value = $airport
$ws.Range("E4").Value = $airport
$ws.Range("E5").Value = $airport
$ws.Range("E6").Value = $gettingAround
$ws.Range("E7").Value = $gettingAround
$ws.Range("E8").Value = $gettingAround
$ws.Range("E9").Value = $gettingAround
$ws.Range("E10").Value = $gettingAround
$ws.Range("E11").Value = $restaurant
$ws.Range("E12").Value = $restaurant
$ws.Range("E13").Value = $restaurant
$ws.Range("E14").Value = $restaurant
$ws.Range("E15").Value = $restaurant

# Match the column widths Excel settled on after adding the new column
# (ColumnWidth is character-width; stored xml width = ColumnWidth + 5/6).
$ws.Cells.Item(1, 2).EntireColumn.ColumnWidth = 13.83
$ws.Cells.Item(1, 3).EntireColumn.ColumnWidth = 63.5
$ws.Cells.Item(1, 4).EntireColumn.ColumnWidth = 38.28
$ws.Cells.Item(1, 5).EntireColumn.ColumnWidth = 16.05

# Restore the selected cell as left by the author
$ws.Range("G19").Select() | Out-Null

# Page setup was touched (portrait orientation) as part of this edit
$ws.PageSetup.Orientation = 1

$wb.Save()
